$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.8775636666666666
$ws.Range("H2").Value = 2.632691
$ws.Range("I2").Value = 0.1887436506618166
$ws.Range("J2").Value = 0.2083714858314108
$ws.Range("M2").Value = 16.14072933333334
$ws.Range("N2").Value = 48.42218800000001
$ws.Range("O2").Value = 0.03423048004954622
$ws.Range("P2").Value = 0.03634868370049611
$ws.Range("Q2").Value = 14.16451761643422
$ws.Range("R2").Value = 127.480658547908
$ws.Range("S2").Value = 0.006460785768457834
$ws.Range("T2").Value = 0.007574029230688357
$ws.Range("G3").Value = 0.8775636666666666
$ws.Range("H3").Value = 2.632691
$ws.Range("I3").Value = 0.1887436506618166
$ws.Range("J3").Value = 0.2083714858314108
$ws.Range("O3").Value = 0.1719151703242873
$ws.Range("P3").Value = 0.1825533892714798
$ws.Range("Q3").Value = 71.13822111363989
$ws.Range("R3").Value = 640.243990022759
$ws.Range("S3").Value = 0.03244789685115398
$ws.Range("T3").Value = 0.03803892096605817
$ws.Range("G4").Value = 0.8775636666666666
$ws.Range("H4").Value = 2.632691
$ws.Range("I4").Value = 0.1887436506618166
$ws.Range("J4").Value = 0.2083714858314108
$ws.Range("M4").Value = 168.70371
$ws.Range("N4").Value = 506.11113
$ws.Range("O4").Value = 0.3577786889414888
$ws.Range("P4").Value = 0.3799182594076638
$ws.Range("Q4").Value = 148.04824632787
$ws.Range("R4").Value = 1332.43421695083
$ws.Range("S4").Value = 0.0675284558798151
$ws.Range("T4").Value = 0.07916413220725826
$ws.Range("G5").Value = 0.8775636666666666
$ws.Range("H5").Value = 2.632691
$ws.Range("I5").Value = 0.1887436506618166
$ws.Range("J5").Value = 0.2083714858314108
$ws.Range("M5").Value = 82.43477250000001
$ws.Range("N5").Value = 164.869545
$ws.Range("O5").Value = 0.1748236883957081
$ws.Range("P5").Value = 0.1237612588479007
$ws.Range("Q5").Value = 72.34176121593251
$ws.Range("R5").Value = 434.050567295595
$ws.Range("S5").Value = 0.0329968611699698
$ws.Range("T5").Value = 0.0257883173945029
$ws.Range("G6").Value = 0.8775636666666666
$ws.Range("H6").Value = 2.632691
$ws.Range("I6").Value = 0.1887436506618166
$ws.Range("J6").Value = 0.2083714858314108
$ws.Range("M6").Value = 123.1883796666667
$ws.Range("N6").Value = 369.565139
$ws.Range("O6").Value = 0.2612519722889696
$ws.Range("P6").Value = 0.2774184087724594
$ws.Range("Q6").Value = 108.1056461510054
$ws.Range("R6").Value = 972.9508153590489
$ws.Range("S6").Value = 0.04930965099241986
$ws.Range("T6").Value = 0.05780608603290305
$ws.Range("G7").Value = 2.458038666666667
$ws.Range("H7").Value = 7.374116000000001
$ws.Range("I7").Value = 0.5286672739959656
$ws.Range("J7").Value = 0.5836444564186148
$ws.Range("M7").Value = 16.14072933333334
$ws.Range("N7").Value = 48.42218800000001
$ws.Range("O7").Value = 0.03423048004954622
$ws.Range("P7").Value = 0.03634868370049611
$ws.Range("Q7").Value = 39.67453680953423
$ws.Range("R7").Value = 357.0708312858081
$ws.Range("S7").Value = 0.01809653457536689
$ws.Range("T7").Value = 0.02121470773990822
$ws.Range("G8").Value = 2.458038666666667
$ws.Range("H8").Value = 7.374116000000001
$ws.Range("I8").Value = 0.5286672739959656
$ws.Range("J8").Value = 0.5836444564186148
$ws.Range("O8").Value = 0.1719151703242873
$ws.Range("P8").Value = 0.1825533892714798
$ws.Range("S8").Value = 0.0908859244538931
$ws.Range("T8").Value = 0.1065462736487287
$ws.Range("G9").Value = 2.458038666666667
$ws.Range("H9").Value = 7.374116000000001
$ws.Range("I9").Value = 0.5286672739959656
$ws.Range("J9").Value = 0.5836444564186148
$ws.Range("M9").Value = 168.70371
$ws.Range("N9").Value = 506.11113
$ws.Range("O9").Value = 0.3577786889414888
$ws.Range("P9").Value = 0.3799182594076638
$ws.Range("Q9").Value = 414.68024239012
$ws.Range("R9").Value = 3732.12218151108
$ws.Range("S9").Value = 0.1891458841765474
$ws.Range("T9").Value = 0.2217371859954922
$ws.Range("G10").Value = 2.458038666666667
$ws.Range("H10").Value = 7.374116000000001
$ws.Range("I10").Value = 0.5286672739959656
$ws.Range("J10").Value = 0.5836444564186148
$ws.Range("M10").Value = 82.43477250000001
$ws.Range("N10").Value = 164.869545
$ws.Range("O10").Value = 0.1748236883957081
$ws.Range("P10").Value = 0.1237612588479007
$ws.Range("Q10").Value = 202.6278582828701
$ws.Range("R10").Value = 1215.76714969722
$ws.Range("S10").Value = 0.09242356277407909
$ws.Range("T10").Value = 0.07223257264596648
$ws.Range("G11").Value = 2.458038666666667
$ws.Range("H11").Value = 7.374116000000001
$ws.Range("I11").Value = 0.5286672739959656
$ws.Range("J11").Value = 0.5836444564186148
$ws.Range("M11").Value = 123.1883796666667
$ws.Range("N11").Value = 369.565139
$ws.Range("O11").Value = 0.2612519722889696
$ws.Range("P11").Value = 0.2774184087724594
$ws.Range("Q11").Value = 302.8018005046804
$ws.Range("R11").Value = 2725.216204542124
$ws.Range("S11").Value = 0.1381153680160791
$ws.Range("T11").Value = 0.1619137163885192
$ws.Range("G12").Value = 1.313898
$ws.Range("H12").Value = 2.627796
$ws.Range("I12").Value = 0.2825890753422177
$ws.Range("J12").Value = 0.2079840577499744
$ws.Range("M12").Value = 16.14072933333334
$ws.Range("N12").Value = 48.42218800000001
$ws.Range("O12").Value = 0.03423048004954622
$ws.Range("P12").Value = 0.03634868370049611
$ws.Range("Q12").Value = 21.207271989608
$ws.Range("R12").Value = 127.243631937648
$ws.Range("S12").Value = 0.009673159705721497
$ws.Range("T12").Value = 0.007559946729899537
$ws.Range("G13").Value = 1.313898
$ws.Range("H13").Value = 2.627796
$ws.Range("I13").Value = 0.2825890753422177
$ws.Range("J13").Value = 0.2079840577499744
$ws.Range("O13").Value = 0.1719151703242873
$ws.Range("P13").Value = 0.1825533892714798
$ws.Range("Q13").Value = 106.508929203734
$ws.Range("R13").Value = 639.053575222404
$ws.Range("S13").Value = 0.04858134901924022
$ws.Range("T13").Value = 0.03796819465669302
$ws.Range("G14").Value = 1.313898
$ws.Range("H14").Value = 2.627796
$ws.Range("I14").Value = 0.2825890753422177
$ws.Range("J14").Value = 0.2079840577499744
$ws.Range("M14").Value = 168.70371
$ws.Range("N14").Value = 506.11113
$ws.Range("O14").Value = 0.3577786889414888
$ws.Range("P14").Value = 0.3799182594076638
$ws.Range("Q14").Value = 221.65946716158
$ws.Range("R14").Value = 1329.95680296948
$ws.Range("S14").Value = 0.1011043488851262
$ws.Range("T14").Value = 0.07901694120491332
$ws.Range("G15").Value = 1.313898
$ws.Range("H15").Value = 2.627796
$ws.Range("I15").Value = 0.2825890753422177
$ws.Range("J15").Value = 0.2079840577499744
$ws.Range("M15").Value = 82.43477250000001
$ws.Range("N15").Value = 164.869545
$ws.Range("O15").Value = 0.1748236883957081
$ws.Range("P15").Value = 0.1237612588479007
$ws.Range("Q15").Value = 108.310882718205
$ws.Range("R15").Value = 433.24353087282
$ws.Range("S15").Value = 0.04940326445165914
$ws.Range("T15").Value = 0.02574036880743131
$ws.Range("G16").Value = 1.313898
$ws.Range("H16").Value = 2.627796
$ws.Range("I16").Value = 0.2825890753422177
$ws.Range("J16").Value = 0.2079840577499744
$ws.Range("M16").Value = 123.1883796666667
$ws.Range("N16").Value = 369.565139
$ws.Range("O16").Value = 0.2612519722889696
$ws.Range("P16").Value = 0.2774184087724594
$ws.Range("Q16").Value = 161.856965667274
$ws.Range("R16").Value = 971.141794003644
$ws.Range("S16").Value = 0.07382695328047059
$ws.Range("T16").Value = 0.05769860635103721
